# Fruta / hortaliza, semanal
# Update the weekly Coco price records: each row's Fecha (D), Volumen (M),
# Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P) and
# Precio $/Kg (S) are refreshed with the latest weekly figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(Fecha, Volumen, PrecioMinimo, PrecioMaximo, PrecioPromedio, PrecioPorKg)
$rows = @{
    2  = @(44893, 80,  21000, 22000, 21625, 1081)
    3  = @(44798, 80,  21000, 22000, 21500, 1075)
    4  = @(44357, 100, 14000, 15000, 14500, 725)
    5  = @(44320, 80,  16000, 17000, 16500, 825)
    6  = @(44533, 100, 16000, 17000, 16500, 825)
    7  = @(44708, 80,  20000, 21000, 20500, 1025)
    8  = @(44890, 80,  20000, 23000, 22250, 1112)
    9  = @(44761, 100, 20000, 21000, 20500, 1025)
    10 = @(45092, 150, 24000, 25000, 24333, 1217)
    11 = @(44792, 100, 21000, 22000, 21500, 1075)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]   # S - Precio $/Kg
}
